# "Add files via upload" — the re-uploaded workbook renamed the
# "ESTADOCIVIL_*" header columns to "DESCRIÇÃO ESTADO CIVIL_*"
# (same suffixes, just a clearer/longer label), leaving every other
# header and value in the sheet untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldPrefix = "ESTADOCIVIL_"
$newPrefix = "DESCRIÇÃO ESTADO CIVIL_"

# Header row lives on row 1; scan its used range and rename any cell
# whose text still starts with the old prefix.
$usedRange = $ws.UsedRange
$lastCol = $usedRange.Columns.Count

for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $text = $cell.Value2
    if ($text -like "$oldPrefix*") {
        $suffix = $text.Substring($oldPrefix.Length)
        $cell.Value = $newPrefix + $suffix
    }
}
